# Edit: 
#  1) Slide 16's table switches to table style {8E33CB3D-23E2-4396-8AFA-14FE86B5AA9E}
#  2) The theme used by the slide master (ppt/theme/theme2.xml, currently "Integral")
#     is swapped with the theme used by the notes master (ppt/theme/theme1.xml,
#     currently "Office Theme") -- i.e. the deck's visible theme reverts to the
#     stock Office color scheme, and the notes-master theme becomes the Integral
#     color scheme.

$p = $ppt.ActivePresentation

# --- 1) Table style id on slide 16 ---------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{8E33CB3D-23E2-4396-8AFA-14FE86B5AA9E}")
    }
}

# --- 2) Swap the two theme color schemes ---------------------------------
# Slide.ThemeColorScheme writes into the theme part used by the slide master
# (ppt/theme/theme2.xml). Slide.NotesPage.ThemeColorScheme writes into the
# theme part used by the notes master (ppt/theme/theme1.xml). Re-pointing
# each of the 12 colour slots gives every run-time consumer of the theme
# (slides + notes pages) the colours from the opposite theme, matching the
# file-content swap in the diff.

$anySlide = $p.Slides.Item(1)

# Office Theme colours (was in theme1.xml) now go to theme2.xml
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)
# Integral colours (was in theme2.xml) now go to theme1.xml
$integralColors = @(0, 16777215, 5332805, 13754083, 3722137, 3646819, 2412774, 38860, 13611854, 10915127, 2465643, 158642)

$slideScheme = $anySlide.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $slideScheme.Item($i).RGB = $officeColors[$i - 1]
}

$notesScheme = $anySlide.NotesPage.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Item($i).RGB = $integralColors[$i - 1]
}

Write-Output "theme colours + table style updated"
